$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Create the "Berechnung" (Calculation) cell style: bold orange font,
#    light-grey solid fill, thin grey border all around.
$berechnung = $wb.Styles.Add("Berechnung")
$berechnung.Font.Color = 32250      # RGB(0x00,0x7D,0xFA) VBA BGR == FA7D00
$berechnung.Font.Bold = $true
$berechnung.Interior.Pattern = 1    # xlSolid
$berechnung.Interior.Color = 15921906  # F2F2F2

# 2) New header row for the second experiment table (G1:K1) plus the
#    "updates per year" column header (M1).
$ws.Range("G1").Value = "num_frames"
$ws.Range("H1").Value = "num_steps"
$ws.Range("I1").Value = "num_processes"
$ws.Range("J1").Value = "num_updates"
$ws.Range("K1").Value = "Bemerkung"
$ws.Range("M1").Value = "updates per year"

# 3) New experiment data / calculation row.
$ws.Range("G2").Value = 2912000
$ws.Range("H2").Value = 14
$ws.Range("I2").Value = 15
$ws.Range("J2").Formula = "=G2/H2/I2"
$ws.Range("M2").Formula = "=364/H2"

# 4) Extra data row (10000 / 14 / 15) appended to the original table.
$ws.Range("A10").Value = 10000
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 15

# 5) Turn column D's per-row formulas into one shared formula D2:D10.
$ws.Range("D2:D10").Formula = "=A2*B2*C2"

# 6) Apply the "Berechnung" style + thin grey border to the highlighted
#    result cells (D1:D10 and J1:J2).
$styledRanges = @("D1", "D2:D10", "J1", "J2")
foreach ($addr in $styledRanges) {
    $rng = $ws.Range($addr)
    $rng.Style = "Berechnung"
    $rng.BorderAround(1, 2, 0, 8355711) | Out-Null
}

# 7) Match the selection shown in the saved workbook.
$ws.Range("I13").Select() | Out-Null

Write-Output "edit applied"
